# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" sheet (fund holdings detail) positioned right
# after "总计" and before the existing "2022-Q2" sheet, and records its
# summary numbers (持有数量(只)=4, 持有市值(亿元)=1.3) as a new row in the
# "总计" sheet (pushing the 2022-Q2 / 2021-Q1 rows down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet.
#    Read the existing rows 2 & 3 first (NOTE: use .Value2, not .Value -
#    this host mis-stringifies bare `.Value` reads), then rewrite rows
#    2-4 top to bottom so the new 2022-Q3 entry lands on row 2.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$oldRow2B = $summary.Range("B2").Value2
$oldRow2C = $summary.Range("C2").Value2
$oldRow2D = $summary.Range("D2").Value2

$oldRow3B = $summary.Range("B3").Value2
$oldRow3C = $summary.Range("C3").Value2
$oldRow3D = $summary.Range("D3").Value2

# Row 4 (was row 3: 2021-Q1) - a brand new row, so give column A the same
# bold / bordered / centered style the other index cells (A2, A3) use.
$summary.Range("A4").Value2 = 2
$a4 = $summary.Range("A4")
$a4.Font.Bold = $true
$a4.Borders.LineStyle = 1
$a4.HorizontalAlignment = -4108
$a4.VerticalAlignment = -4160
$summary.Range("B4").Value2 = $oldRow3B
$summary.Range("C4").Value2 = $oldRow3C
$summary.Range("D4").Value2 = $oldRow3D

# Row 3 (was row 2: 2022-Q2) - cell already styled, only values move.
$summary.Range("A3").Value2 = 1
$summary.Range("B3").Value2 = $oldRow2B
$summary.Range("C3").Value2 = $oldRow2C
$summary.Range("D3").Value2 = $oldRow2D

# Row 2 (new: 2022-Q3) - cell already styled, only values change.
$summary.Range("A2").Value2 = 0
$summary.Range("B2").Value2 = "2022-Q3"
$summary.Range("C2").Value2 = 4
$summary.Range("D2").Value2 = 1.3

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet before "2022-Q2" and populate it
#    with the fund holdings detail.
# ---------------------------------------------------------------------
$existingQ2 = $wb.Worksheets.Item("2022-Q2")
$ws = $wb.Worksheets.Add($existingQ2)
$ws.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 2)
    $cell.Value2 = $headers[$i]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Columns B-G are text (fund code / percentages stored as strings, not
# numbers - e.g. "010695" and "0.6630" must keep their leading/trailing
# zeros), so they are entered with a leading apostrophe to force text
# storage. Columns A and H are real numbers.
$data = @(
    @(0, "'010695", "华夏磐益一年定期开放混合", "'15.90", "'99.95", "'4.17", "'0.6630", 8),
    @(1, "'009837", "华夏磐锐一年定期开放混合A", "'14.02", "'94.15", "'4.20", "'0.5888", 8),
    @(2, "'620001", "金元顺安宝石动力混合", "'1.01", "'40.12", "'3.45", "'0.0348", 4),
    @(3, "'009838", "华夏磐锐一年定期开放混合C", "'0.39", "'94.15", "'4.20", "'0.0164", 8)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    $rowNum = $r + 2

    $aCell = $ws.Cells.Item($rowNum, 1)
    $aCell.Value2 = $row[0]
    $aCell.Font.Bold = $true
    $aCell.Borders.LineStyle = 1
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160

    $ws.Cells.Item($rowNum, 2).Value2 = $row[1]
    $ws.Cells.Item($rowNum, 3).Value2 = $row[2]
    $ws.Cells.Item($rowNum, 4).Value2 = $row[3]
    $ws.Cells.Item($rowNum, 5).Value2 = $row[4]
    $ws.Cells.Item($rowNum, 6).Value2 = $row[5]
    $ws.Cells.Item($rowNum, 7).Value2 = $row[6]
    $ws.Cells.Item($rowNum, 8).Value2 = $row[7]
}

$ws.Range("A1").Select()
